$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (Volume number + week-covering date range) ---
$ws.Range("A8").Value = "Volume 31   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"

# --- Column E width (auto bestFit widened due to longer "***.* "-style values) ---
$ws.Columns("E").ColumnWidth = 7.433768

# --- Data cell updates ---
# Row 14
$ws.Range("L14").Value = -41.176470588235
$ws.Range("M14").Value = -54.545454545454
$ws.Range("N14").Value = -77.011494252873
# Row 15
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 18.181818181818
$ws.Range("I15").Value = 114
$ws.Range("J15").Value = 91
$ws.Range("K15").Value = 25.274725274725
$ws.Range("L15").Value = -5
$ws.Range("M15").Value = 29.545454545454
$ws.Range("N15").Value = -48.878923766816
# Row 16
$ws.Range("C16").Value = 32
$ws.Range("D16").Value = 28
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 120
$ws.Range("G16").Value = 116
$ws.Range("H16").Value = 3.448275862068
$ws.Range("I16").Value = 877
$ws.Range("J16").Value = 838
$ws.Range("K16").Value = 4.653937947494
$ws.Range("L16").Value = -6.203208556149
$ws.Range("M16").Value = -33.257229832572
$ws.Range("N16").Value = -81.594963273872
# Row 17
$ws.Range("C17").Value = 54
$ws.Range("D17").Value = 55
$ws.Range("E17").Value = -1.818181818181
$ws.Range("F17").Value = 226
$ws.Range("G17").Value = 236
$ws.Range("H17").Value = -4.237288135593
$ws.Range("I17").Value = 1974
$ws.Range("J17").Value = 1785
$ws.Range("K17").Value = 10.588235294117
$ws.Range("L17").Value = 13.513513513513
$ws.Range("M17").Value = 83.286908077994
$ws.Range("N17").Value = -15.387912558937
# Row 18
$ws.Range("C18").Value = 31
$ws.Range("D18").Value = 18
$ws.Range("E18").Value = 72.222222222222
$ws.Range("F18").Value = 92
$ws.Range("G18").Value = 99
$ws.Range("H18").Value = -7.070707070707
$ws.Range("I18").Value = 655
$ws.Range("J18").Value = 732
$ws.Range("K18").Value = -10.51912568306
$ws.Range("L18").Value = -7.616361071932
$ws.Range("M18").Value = -49.381761978361
$ws.Range("N18").Value = -88.606714211167
# Row 19
$ws.Range("C19").Value = 67
$ws.Range("D19").Value = 83
$ws.Range("E19").Value = -19.277108433734
$ws.Range("F19").Value = 274
$ws.Range("G19").Value = 340
$ws.Range("H19").Value = -19.411764705882
$ws.Range("I19").Value = 2039
$ws.Range("J19").Value = 2290
$ws.Range("K19").Value = -10.960698689956
$ws.Range("L19").Value = -13.820794590025
$ws.Range("M19").Value = 18.271461716937
$ws.Range("N19").Value = -60.637065637065
# Row 20
$ws.Range("C20").Value = 44
$ws.Range("D20").Value = 43
$ws.Range("E20").Value = 2.325581395348
$ws.Range("F20").Value = 158
$ws.Range("G20").Value = 166
$ws.Range("H20").Value = -4.819277108433
$ws.Range("I20").Value = 1217
$ws.Range("J20").Value = 1120
$ws.Range("K20").Value = 8.660714285714
$ws.Range("L20").Value = 24.948665297741
$ws.Range("M20").Value = 13.844714686623
$ws.Range("N20").Value = -89.752441899629
# Row 21
$ws.Range("C21").Value = 234
$ws.Range("E21").Value = 1.298701298701
$ws.Range("F21").Value = 884
$ws.Range("G21").Value = 968
$ws.Range("H21").Value = -8.677685950413
$ws.Range("I21").Value = 6896
$ws.Range("J21").Value = 6866
$ws.Range("K21").Value = 0.436935624817
$ws.Range("L21").Value = 0.27628326305
$ws.Range("M21").Value = 4.326777609682
$ws.Range("N21").Value = -77.175388077979
# Row 22
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 13
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = 8.333333333333
$ws.Range("I22").Value = 82
$ws.Range("J22").Value = 76
$ws.Range("K22").Value = 7.894736842105
$ws.Range("L22").Value = 12.328767123287
$ws.Range("M22").Value = 3.79746835443
# Row 23
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -28.571428571428
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = -4.761904761904
$ws.Range("I23").Value = 148
$ws.Range("J23").Value = 156
$ws.Range("K23").Value = -5.128205128205
$ws.Range("L23").Value = 2.068965517241
$ws.Range("M23").Value = 59.139784946236
# Row 24
$ws.Range("C24").Value = 208
$ws.Range("D24").Value = 172
$ws.Range("E24").Value = 20.930232558139
$ws.Range("F24").Value = 789
$ws.Range("G24").Value = 693
$ws.Range("H24").Value = 13.852813852813
$ws.Range("I24").Value = 5767
$ws.Range("J24").Value = 5848
$ws.Range("K24").Value = -1.385088919288
$ws.Range("L24").Value = -7.742761158214
$ws.Range("M24").Value = 48.710675605982
# Row 25
$ws.Range("C25").Value = 86
$ws.Range("D25").Value = 56
$ws.Range("E25").Value = 53.571428571428
$ws.Range("F25").Value = 340
$ws.Range("G25").Value = 241
$ws.Range("H25").Value = 41.078838174273
$ws.Range("I25").Value = 2460
$ws.Range("J25").Value = 1987
$ws.Range("K25").Value = 23.804730749874
$ws.Range("L25").Value = 10.71107110711
# Row 26
$ws.Range("C26").Value = 124
$ws.Range("D26").Value = 77
$ws.Range("E26").Value = 61.038961038961
$ws.Range("F26").Value = 433
$ws.Range("G26").Value = 377
$ws.Range("H26").Value = 14.854111405835
$ws.Range("I26").Value = 3257
$ws.Range("J26").Value = 2847
$ws.Range("K26").Value = 14.401123990165
$ws.Range("L26").Value = 26.830218068535
$ws.Range("M26").Value = 10.857726344452
# Row 27
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 75
$ws.Range("F27").Value = 25
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = 78.571428571428
$ws.Range("I27").Value = 181
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 13.125
$ws.Range("L27").Value = -5.729166666666
# Row 28
$ws.Range("F28").Value = 47
$ws.Range("G28").Value = 32
$ws.Range("H28").Value = 46.875
$ws.Range("I28").Value = 285
$ws.Range("J28").Value = 268
$ws.Range("K28").Value = 6.343283582089
$ws.Range("L28").Value = 2.888086642599
# Row 29
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 66.666666666666
$ws.Range("I29").Value = 72
$ws.Range("J29").Value = 59
$ws.Range("K29").Value = 22.033898305084
$ws.Range("L29").Value = -41.935483870967
$ws.Range("M29").Value = -43.307086614173
$ws.Range("N29").Value = -76.699029126213
# Row 30
$ws.Range("C30").Value = 1
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 33.333333333333
$ws.Range("I30").Value = 59
$ws.Range("J30").Value = 43
$ws.Range("K30").Value = 37.209302325581
$ws.Range("L30").Value = -37.234042553191
$ws.Range("M30").Value = -42.156862745098
$ws.Range("N30").Value = -79.225352112676
# Row 31
$ws.Range("I31").Value = 22
$ws.Range("K31").Value = -42.105263157894
$ws.Range("L31").Value = -12
# Row 33
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = -75
$ws.Range("J33").Value = 27
$ws.Range("K33").Value = -51.851851851851

# --- Type-changing cells (string <-> number), with style fix-ups to match target formatting ---

# D29: style 14 -> 15, type s -> n
$ws.Range("D29").Value = 1
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)

# E29: style 14 -> 16, type s -> n
$ws.Range("E29").Value = 100
$ws.Range("H29").Copy()
$ws.Range("E29").PasteSpecial(-4122)

# D30: style 14 -> 15, type s -> n
$ws.Range("D30").Value = 1
$ws.Range("C29").Copy()
$ws.Range("D30").PasteSpecial(-4122)

# E30: style 14 -> 16, type s -> n
$ws.Range("E30").Value = 0
$ws.Range("H29").Copy()
$ws.Range("E30").PasteSpecial(-4122)

# C31: style 14 -> 15, type s -> n
$ws.Range("C31").Value = 1
$ws.Range("C29").Copy()
$ws.Range("C31").PasteSpecial(-4122)

# C33: style 15 -> 14, type n -> s
$ws.Range("C33").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C33").PasteSpecial(-4122)

$excel.CutCopyMode = $false
$ws.Calculate()
